$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New last row for this sheet's table (row 34 is the current last data row).
$srcRow = 34
$dstRow = 35

# Copy the previous row's cell formats (and values, as a base) down to the
# new row so the new row picks up the same styles (bold/border on col A,
# date-time number format on col E) without creating new style entries.
$srcRange = "A" + $srcRow + ":V" + $srcRow
$dstRange = "A" + $dstRow + ":V" + $dstRow
$ws.Range($srcRange).Copy($ws.Range($dstRange))

# Now overwrite with the actual values for the new match record.
$ws.Cells.Item($dstRow, 1).Value = 34
$ws.Cells.Item($dstRow, 2).Value = "india"
$ws.Cells.Item($dstRow, 3).Value = "isl"
$ws.Cells.Item($dstRow, 4).Value = "2023-2024"
$ws.Cells.Item($dstRow, 5).Value = 45237.64583333334
$ws.Cells.Item($dstRow, 6).Value = "Punjab"
$ws.Cells.Item($dstRow, 7).Value = 1
$ws.Cells.Item($dstRow, 8).Value = "Hyderabad"
$ws.Cells.Item($dstRow, 9).Value = 1
$ws.Cells.Item($dstRow, 10).Value = 3.83
$ws.Cells.Item($dstRow, 11).Value = "06/11/2023 15:18"
$ws.Cells.Item($dstRow, 12).Value = 3.37
$ws.Cells.Item($dstRow, 13).Value = "07/11/2023 15:22"
$ws.Cells.Item($dstRow, 14).Value = 3.49
$ws.Cells.Item($dstRow, 15).Value = "06/11/2023 15:18"
$ws.Cells.Item($dstRow, 16).Value = 3.65
$ws.Cells.Item($dstRow, 17).Value = "07/11/2023 15:22"
$ws.Cells.Item($dstRow, 18).Value = 1.9
$ws.Cells.Item($dstRow, 19).Value = "06/11/2023 15:18"
$ws.Cells.Item($dstRow, 20).Value = 2.1
$ws.Cells.Item($dstRow, 21).Value = "07/11/2023 15:22"
$ws.Cells.Item($dstRow, 22).Value = "https://www.betexplorer.com/football/india/isl/minerva-punjab-hyderabad/bLmEZm0S/"

Write-Host "Row 35 added"
